# Update countries & provincias Spain
# Daily COVID data refresh: San Marino and Gabon rows moved up (re-sorted by
# "Casos totales" descending) with refreshed figures; the rows that used to
# sit above them shift down one slot (their own values unchanged); a handful
# of other countries' daily figures were refreshed too; the "datos
# actualizados" timestamp moved from 21:22 to 21:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Title / timestamp row -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 21:52"

# --- Simple daily figure refreshes (no row movement) ------------------------
# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 457101
$ws.Cells.Item(4, 3).Value = 22174
$ws.Cells.Item(4, 4).Value = 24910
$ws.Cells.Item(4, 5).Value = 415945
$ws.Cells.Item(4, 7).Value = 1458
$ws.Cells.Item(4, 8).Value = 16246

# Alemania (row 8)
$ws.Cells.Item(8, 4).Value = 50557
$ws.Cells.Item(8, 5).Value = 62515

# Canada (row 16)
$ws.Cells.Item(16, 2).Value = 20696
$ws.Cells.Item(16, 3).Value = 1258
$ws.Cells.Item(16, 4).Value = 5206
$ws.Cells.Item(16, 5).Value = 14987

# Noruega (row 26)
$ws.Cells.Item(26, 2).Value = 6162
$ws.Cells.Item(26, 3).Value = 120
$ws.Cells.Item(26, 5).Value = 6022

# --- San Marino jumps above Banglades (rows 101-104) ------------------------
# New order: San Marino (refreshed), Banglades, Mauricio, Ghana (each of the
# latter three keep their previous values, just shifted down one row).
Set-Row 101 @("San Marino", 333, 25, 49, 250, 14, 0, 34)
Set-Row 102 @("Banglades", 330, 112, 33, 276, 1, 1, 21)
Set-Row 103 @("Mauricio", 314, 41, 23, 284, 3, 0, 7)
Set-Row 104 @("Ghana", 313, 0, 34, 273, 2, 0, 6)

# --- Gabon jumps above San Martin (Parte Holandesa) (rows 147-154) ----------
Set-Row 147 @("Gabon", 44, 10, 1, 42, 0, 0, 1)
Set-Row 148 @("San Martin (Parte Holandesa)", 43, 3, 1, 36, 2, 0, 6)
Set-Row 149 @("Bahamas", 40, 0, 5, 28, 1, 0, 7)
Set-Row 150 @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
Set-Row 151 @("Zambia", 39, 0, 24, 14, 1, 0, 1)
Set-Row 152 @("Bermudas", 39, 0, 23, 13, 0, 0, 3)
Set-Row 153 @("Guyana", 37, 0, 8, 23, 4, 0, 6)
Set-Row 154 @("Guinea-Bisau", 36, 3, 0, 36, 0, 0, 0)
